$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text helper cells (column D numeric-looking values) use a temporary
# Text number-format so Excel keeps the literal string (incl. trailing zeros
# and does not coerce to a float), then ClearFormats() removes the temporary
# style again so the cell keeps its original (default) style index.

$ws.Range("D2").Value = '46.729.02'
$ws.Range("E2").Value = '  -0.75%  '

$ws.Range("D3").Value = '2.256.77'
$ws.Range("E3").Value = '  -4.05%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '297.55'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.72'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.36%  '

$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.500'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -7.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.50'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.22%  '

$ws.Range("E12").Value = '  -6.33%  '

$ws.Range("D14").Value = '2.601.06'
$ws.Range("E14").Value = '  -3.93%  '

$ws.Range("D15").Value = '2.261.51'
$ws.Range("E15").Value = '  -3.76%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.51'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.25%  '

$ws.Range("D17").Value = '46.694.46'
$ws.Range("E17").Value = '  -0.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.789'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -5.10%  '

$ws.Range("D19").Value = '0.0₃0965'
$ws.Range("E19").Value = '  +1.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.28'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -10.40%  '

$ws.Range("E21").Value = '  -7.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.55'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.60'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.26%  '

$ws.Range("E24").Value = '  -7.35%  '

$ws.Range("E25").Value = '  +0.28%  '

$ws.Range("E26").Value = '  -7.77%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '40.94'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.40%  '

$ws.Range("E28").Value = '  -3.58%  '

$ws.Range("E29").Value = '  -4.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.96'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.22%  '

$ws.Range("E31").Value = '  +6.84%  '

$ws.Range("E32").Value = '  +3.81%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '143.66'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.27'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -8.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0762'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.110'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.19%  '

$ws.Range("E37").Value = '  -3.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.13'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +9.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.64'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -10.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.77'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0293'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -7.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.06'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -10.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '92.09'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +13.39%  '

$ws.Range("D45").Value = '1.772.94'
$ws.Range("E45").Value = '  -4.32%  '

$ws.Range("E46").Value = '  -7.89%  '

$ws.Range("B47").Value = 'ordi'
$ws.Range("C47").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '69.65'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.45%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.182'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -7.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.74'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.21%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.79'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.98%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '93.43'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.34%  '
